$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.667.17"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "'1.584.73"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'206.43"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'22.30"
$ws.Range("E8").Value = "  -4.81%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'1.809.32"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "'1.574.80"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "'3.87"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("D16").Value = "'27.636.42"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'63.07"
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("D18").Value = "'218.41"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "'0.0" + [char]8323 + "0694"
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("E20").Value = "  -4.93%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("D23").Value = "'9.53"
$ws.Range("E23").Value = "  -5.41%  "
$ws.Range("D24").Value = "'1.98"
$ws.Range("E24").Value = "  -4.97%  "
$ws.Range("D25").Value = "'153.55"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'6.72"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").Value = "'15.07"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "'1.15"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "'0.0465"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("D33").Value = "'1.378.81"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").Value = "'1.52"
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").Value = "'0.540"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("D40").Value = "'0.821"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").Value = "'63.86"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").Value = "'5.23"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").Value = "'1.719.86"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").Value = "'87.66"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'0.0" + [char]8326 + "0101"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "'0.0975"
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("E51").Value = "  -1.30%  "
